# Update gh-pages to output generated at 456a3b4
#
# This refreshes the "F" column (想去人数 / interest count) values on the
# 展览 (Worksheets index 1), 演出 (index 2) and 全部类型 (index 4) sheets
# to a newer data snapshot.

$wb = $excel.ActiveWorkbook

# 展览 (sheet1)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F6").Value  = 1368
$ws1.Range("F9").Value  = 405
$ws1.Range("F10").Value = 473
$ws1.Range("F11").Value = 837
$ws1.Range("F12").Value = 539
$ws1.Range("F17").Value = 1070
$ws1.Range("F18").Value = 521
$ws1.Range("F19").Value = 302
$ws1.Range("F20").Value = 432
$ws1.Range("F22").Value = 250
$ws1.Range("F23").Value = 33
$ws1.Range("F25").Value = 502
$ws1.Range("F27").Value = 2
$ws1.Range("F28").Value = 357

# 演出 (sheet2)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F3").Value  = 388
$ws2.Range("F13").Value = 10

# 全部类型 (sheet4)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F7").Value  = 1368
$ws4.Range("F10").Value = 388
$ws4.Range("F13").Value = 405
$ws4.Range("F16").Value = 473
$ws4.Range("F17").Value = 837
$ws4.Range("F18").Value = 539
$ws4.Range("F23").Value = 1070
$ws4.Range("F24").Value = 521
$ws4.Range("F27").Value = 302
$ws4.Range("F28").Value = 432
$ws4.Range("F32").Value = 250
$ws4.Range("F33").Value = 33
$ws4.Range("F37").Value = 502
$ws4.Range("F38").Value = 10
$ws4.Range("F41").Value = 2
$ws4.Range("F42").Value = 357

$wb.Save()
